# NewEventTest.xlsx - "added excel config 2"
# Reworks the sample event row: new From/To date values (with an updated
# date/time display format), a numeric-looking Title, a lower-cased test
# value, new decimal sample values and new symbol sample values - plus the
# header renames and column layout tweaks that go with them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Header row tweaks -----------------------------------------------
$ws.Range("B1").Value = "From Date"
$ws.Range("C1").Value = "To Date"

# --- Date/time values (B2:C2) -----------------------------------------
# Give the From/To date columns a two-digit-hour display format, then
# enter the new sample values. A leading apostrophe keeps Excel from
# re-parsing the already-formatted date strings as a live date serial
# (matches the quotePrefix behaviour of the original cells).
$ws.Range("B2:C2").NumberFormat = "dd\-mmm\-yyyy\ hh:mm"
$ws.Range("B2").Value = "'02-Apr-2018  10:20"
$ws.Range("C2").Value = "'05-Apr-2018  10:21"

# --- Title (A2) now holds a numeric-looking string --------------------
$ws.Range("A2").Value = "'123"

# --- Remnote (K2): lower-case test value -------------------------------
$ws.Range("K2").Value = "test"

# --- New decimal sample values (N2:R2), kept as text -------------------
$ws.Range("N2").Value = "'1.1"
$ws.Range("O2").Value = "'2.1"
$ws.Range("P2").Value = "'3.1"
$ws.Range("Q2").Value = "'4.4"
$ws.Range("R2").Value = "'5.1"

# --- New symbol sample values (S2:V2) -----------------------------------
$ws.Range("S2").Value = "##"
$ws.Range("T2").Value = "$$"
$ws.Range("U2").Value = "!!"
$ws.Range("V2").Value = "Test11$"

# --- Column layout: widen A (now bestfit-worthy), resize B/C -----------
$ws.Columns("A").ColumnWidth = 12.25
$ws.Columns("B").ColumnWidth = 16.25
$ws.Columns("C").ColumnWidth = 16.59

# --- Active cell moves to C2 --------------------------------------------
$ws.Range("C2").Select()
